$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row above the current row 35 (shifts existing rows 35-61 down to 36-62).
$ws.Rows.Item(35).Insert()

# Populate the newly inserted row 35 with the new weekly price record.
$ws.Range("A35").Value = 1
$ws.Range("B35").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C35").Value = 'Arica y Parinacota'
$ws.Range("D35").Value = 44658
$ws.Range("E35").Value = 15
$ws.Range("F35").Value = 100112009
$ws.Range("G35").Value = 'Acelga'
$ws.Range("H35").Value = 'Sin especificar'
$ws.Range("I35").Value = 'Primera'
$ws.Range("J35").Value = 200
$ws.Range("K35").Value = 2500
$ws.Range("L35").Value = 2800
$ws.Range("M35").Value = 2650
$ws.Range("N35").Value = '$/atado 2,5 a 3 kilos'
$ws.Range("O35").Value = 'Región de Arica y Parinacota'
$ws.Range("P35").Value = 883
$ws.Range("Q35").Value = 3
$ws.Range("R35").Value = 'Hortaliza'
